$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$bVal = 8.8839854167789998
$cVal = 8.9292780144256092
$dVal = 99.781580779108495

for ($i = 202; $i -le 251; $i++) {
    $aVal = $i - 1
    $ws.Cells.Item($i, 1).Value = $aVal
    $ws.Cells.Item($i, 2).Value = $bVal
    $ws.Cells.Item($i, 3).Value = $cVal
    $ws.Cells.Item($i, 4).Value = $dVal
}

$ws.Range("B201:D251").Select()
$excel.ActiveWindow.ScrollRow = 198
